# Auto-generated Excel COM-interop script
# Applies data refresh updates to the Aegis_Profits workbook (columns H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 11499.75
$ws.Range("I47").Value = 4000
$ws.Range("J47").Value = 13999.667
$ws.Range("K47").Value = 4000
$ws.Range("L47").Value = 13999.667
$ws.Range("M47").Value = -3028
$ws.Range("N47").Value = -15943.667
$ws.Range("H96").Value = 1339.4333
$ws.Range("I96").Value = 433.3684
$ws.Range("J96").Value = 2904.4546
$ws.Range("K96").Value = 1300.1052
$ws.Range("L96").Value = 8713.363799999999
$ws.Range("M96").Value = 72.89480000000003
$ws.Range("N96").Value = -11459.3638
$ws.Range("H106").Value = 1975.9
$ws.Range("I106").Value = 1695.4445
$ws.Range("K106").Value = 1695.4445
$ws.Range("M106").Value = -1064.4445
$ws.Range("H135").Value = 4013.0625
$ws.Range("I135").Value = 1192.75
$ws.Range("J135").Value = 5705.25
$ws.Range("K135").Value = 10734.75
$ws.Range("L135").Value = 51347.25
$ws.Range("M135").Value = -8199.75
$ws.Range("N135").Value = -56417.25
$ws.Range("H137").Value = 2058.2632
$ws.Range("I137").Value = 1585.1538
$ws.Range("K137").Value = 4755.4614
$ws.Range("M137").Value = -2205.4614
$ws.Range("H138").Value = 2513.88
$ws.Range("I138").Value = 1180.95
$ws.Range("J138").Value = 2847.1125
$ws.Range("K138").Value = 3542.85
$ws.Range("L138").Value = 8541.337500000001
$ws.Range("M138").Value = 1597.15
$ws.Range("N138").Value = -18821.3375

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1021.4737
$ws.Range("I74").Value = 980.96295
$ws.Range("K74").Value = 980.96295
$ws.Range("M74").Value = -106.96295
$ws.Range("H77").Value = 1021.4737
$ws.Range("I77").Value = 980.96295
$ws.Range("K77").Value = 4904.81475
$ws.Range("M77").Value = -536.8147499999995
$ws.Range("H122").Value = 2272.25
$ws.Range("I122").Value = 1884.8948
$ws.Range("J122").Value = 3090
$ws.Range("K122").Value = 5654.6844
$ws.Range("L122").Value = 9270
$ws.Range("M122").Value = -3204.6844
$ws.Range("N122").Value = -14170

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 56654.85
$ws.Range("I86").Value = 70461.875
$ws.Range("J86").Value = 1426.75
$ws.Range("K86").Value = 70461.875
$ws.Range("L86").Value = 1426.75
$ws.Range("M86").Value = -69338.875
$ws.Range("H89").Value = 56654.85
$ws.Range("I89").Value = 70461.875
$ws.Range("J89").Value = 1426.75
$ws.Range("K89").Value = 352309.375
$ws.Range("L89").Value = 7133.75
$ws.Range("M89").Value = -346693.375
$ws.Range("N86").Value = -3672.75
$ws.Range("N89").Value = -18365.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1108.3334
$ws.Range("I8").Value = 800
$ws.Range("J8").Value = 1262.5
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 1262.5
$ws.Range("M8").Value = -660
$ws.Range("N8").Value = -1542.5
$ws.Range("H51").Value = 7979.6
$ws.Range("J51").Value = 7979.6
$ws.Range("L51").Value = 7979.6
$ws.Range("N51").Value = -9451.6
$ws.Range("H61").Value = 7979.6
$ws.Range("J61").Value = 7979.6
$ws.Range("L61").Value = 7979.6
$ws.Range("N61").Value = -8675.6
$ws.Range("H141").Value = 106650
$ws.Range("I141").Value = 125000
$ws.Range("J141").Value = 69950
$ws.Range("K141").Value = 125000
$ws.Range("L141").Value = 69950
$ws.Range("M141").Value = -119820
$ws.Range("N141").Value = -80310

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1001.1111
$ws.Range("J34").Value = 1270
$ws.Range("L34").Value = 3810
$ws.Range("N34").Value = -3978
$ws.Range("H52").Value = 2421
$ws.Range("J52").Value = 2421
$ws.Range("L52").Value = 7263
$ws.Range("N52").Value = -7795
$ws.Range("H81").Value = 1966.9565
$ws.Range("I81").Value = 1088
$ws.Range("J81").Value = 2211.111
$ws.Range("K81").Value = 3264
$ws.Range("L81").Value = 6633.333
$ws.Range("M81").Value = -2141
$ws.Range("N81").Value = -8879.332999999999
$ws.Range("H84").Value = 1966.9565
$ws.Range("I84").Value = 1088
$ws.Range("J84").Value = 2211.111
$ws.Range("K84").Value = 9792
$ws.Range("L84").Value = 19899.999
$ws.Range("M84").Value = -4176
$ws.Range("N84").Value = -31131.999
$ws.Range("H125").Value = 1247.5
$ws.Range("I125").Value = 1247.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3742.5
$ws.Range("L125").Value = 0
$ws.Range("H131").Value = 1220.89
$ws.Range("I131").Value = 470.05
$ws.Range("J131").Value = 1408.6
$ws.Range("K131").Value = 1410.15
$ws.Range("L131").Value = 4225.799999999999
$ws.Range("M131").Value = 3629.85
$ws.Range("N131").Value = -14305.8
$ws.Range("M125").Value = 1177.5
$ws.Range("N125").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 179.66667
$ws.Range("I2").Value = 67.22221999999999
$ws.Range("J2").Value = 517
$ws.Range("K2").Value = 67.22221999999999
$ws.Range("L2").Value = 517
$ws.Range("M2").Value = 45.77778000000001
$ws.Range("N2").Value = -743
$ws.Range("H43").Value = 2000.7878
$ws.Range("I43").Value = 960.3461
$ws.Range("J43").Value = 5865.2856
$ws.Range("K43").Value = 960.3461
$ws.Range("L43").Value = 5865.2856
$ws.Range("M43").Value = -809.3461
$ws.Range("N43").Value = -6167.2856
$ws.Range("H46").Value = 12453.2
$ws.Range("J46").Value = 12453.2
$ws.Range("L46").Value = 12453.2
$ws.Range("N46").Value = -12765.2
$ws.Range("H57").Value = 11175
$ws.Range("H132").Value = 4119.2104
$ws.Range("I132").Value = 2712.923
$ws.Range("K132").Value = 8138.768999999999
$ws.Range("M132").Value = -5608.768999999999
$ws.Range("H134").Value = 14493.5
$ws.Range("J134").Value = 14493.5
$ws.Range("L134").Value = 43480.5
$ws.Range("N134").Value = -48550.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2487.8572
$ws.Range("I7").Value = 1808.5714
$ws.Range("J7").Value = 3167.1428
$ws.Range("K7").Value = 1808.5714
$ws.Range("L7").Value = 3167.1428
$ws.Range("M7").Value = -1696.5714
$ws.Range("N7").Value = -3391.1428
$ws.Range("H126").Value = 2487.8572
$ws.Range("I126").Value = 1808.5714
$ws.Range("J126").Value = 3167.1428
$ws.Range("K126").Value = 5425.7142
$ws.Range("L126").Value = 9501.428400000001
$ws.Range("M126").Value = -2955.7142
$ws.Range("N126").Value = -14441.4284
$ws.Range("H132").Value = 3703.2964
$ws.Range("I132").Value = 3859.8
$ws.Range("K132").Value = 11579.4
$ws.Range("M132").Value = -9049.400000000001
$ws.Range("H135").Value = 46161.668
$ws.Range("J135").Value = 46161.668
$ws.Range("L135").Value = 46161.668
$ws.Range("N135").Value = -56301.668

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2001.4783
$ws.Range("I126").Value = 1777.4375
$ws.Range("K126").Value = 5332.3125
$ws.Range("M126").Value = -2862.3125
